$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 55
$ws1.Range("F4").Value = 1815
$ws1.Range("F5").Value = 801
$ws1.Range("F6").Value = 429
$ws1.Range("F7").Value = 219

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 55
$ws4.Range("F4").Value = 1815
$ws4.Range("F6").Value = 801
$ws4.Range("F7").Value = 429
$ws4.Range("F8").Value = 219
